# DPLKINV001-018 - Setup Mapping PVR-TB Investasi
# Replace the generic View/Tambah/Ubah/Hapus step-by-step instructions in
# column D (rows 2-5) with short, specific descriptions that reference the
# "Setup Mapping PVR-TB" feature, and refresh the row heights / selection
# that Excel recalculates as a consequence of the shorter text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Tambah" (Add) scenario
$ws.Range("D2").Value = "Tambah Setup Mapping PVR-TB"

# Row 3 - "View" scenario
$ws.Range("D3").Value = "View Setup Mapping PVR-TB"

# Row 4 - "Ubah" (Edit) scenario
$ws.Range("D4").Value = "Ubah Setup Mapping PVR-TB"

# Row 5 - "Hapus" (Delete) scenario
$ws.Range("D5").Value = "Hapus Setup Mapping PVR-TB"

# The shorter text no longer needs as much vertical space when wrapped, so
# the rows shrink from their previous (75/45/60/45) heights.
$ws.Rows(2).RowHeight = 30
$ws.Rows(3).RowHeight = 30
$ws.Rows(4).AutoFit()
$ws.Rows(5).RowHeight = 30

# Reflect the author's final selection / active cell.
$ws.Range("D5").Select() | Out-Null
